# Register user and add product to wishlist
# Adds "Nav menu options" / "Checkout details" sections (rows 7-12) to the
# test-data sheet, and applies wrap-text formatting across the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: clone the green header look from A1 first -----------
$ws.Range("A1").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

# --- New cell values --------------------------------------------------------
# (entered in this particular order so the shared-strings table comes out in
#  the same sequence as it would from interactive data entry)
$ws.Range("A9").Value  = "Computers"
$ws.Range("A7").Value  = "Nav menu options"
$ws.Range("B9").Value  = " Desktops "
$ws.Range("C9").Value  = "Digital Storm VANQUISH 3 Custom Performance PC"
$ws.Range("C8").Value  = "Destop product to add:"
$ws.Range("B8").Value  = "Computes sub category option:"
$ws.Range("A8").Value  = "Nav menu option:"
$ws.Range("A10").Value = "Checkout details"
$ws.Range("A11").Value = "City:"
$ws.Range("A12").Value = "Mohali"
$ws.Range("B11").Value = "Address 1:"
$ws.Range("B12").Value = "Sector 71"

# --- Wrap text: section headers --------------------------------------------
$ws.Range("A1").WrapText  = $true
$ws.Range("A4").WrapText  = $true
$ws.Range("A7").WrapText  = $true
$ws.Range("A10").WrapText = $true

# --- Wrap text: data rows ---------------------------------------------------
$ws.Range("A2:D2").WrapText   = $true
$ws.Range("A3:D3").WrapText   = $true
$ws.Range("A5:B5").WrapText   = $true
$ws.Range("A6:B6").WrapText   = $true
$ws.Range("A8:C8").WrapText   = $true
$ws.Range("A9:C9").WrapText   = $true
$ws.Range("A11:B11").WrapText = $true
$ws.Range("A12:B12").WrapText = $true

# --- Row 9 is taller because of the wrapped long product name --------------
$ws.Rows(9).RowHeight = 28.8

# --- Selection ends on B12, matching the last cell touched ------------------
[void]$ws.Range("B12").Select()

$wb.Save()
